$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing source text (row 2, column B) before shifting rows down
$sourceText = $ws.Range("B2").Value2

# Insert a new row at row 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Fill in the new row's values: year 2008, source = same text as old B2
$ws.Range("A2").Value = 2008
$ws.Range("B2").Value = $sourceText

# Update the selection to reflect the edit location
$ws.Range("B2").Select()
